$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of effort data: 05/10/2012, 1.75h effort, new description text
$ws.Range("A19").Value = 41187

$ws.Range("B19").Value = 1.75

$ws.Range("D19").Value = "Code cleanup check of all test cases as preparation of re-implementation index->pointer "

# Update the selection to reflect where the cursor now sits after the edit
$ws.Range("A20").Select()
